$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4000
$ws.Range("J17").Value = 4000
$ws.Range("L17").Value = 12000
$ws.Range("N17").Value = -12336

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3225.2856
$ws.Range("J76").Value = 3942
$ws.Range("L76").Value = 3942
$ws.Range("N76").Value = -4572

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3225.2856
$ws.Range("J79").Value = 3942
$ws.Range("L79").Value = 3942
$ws.Range("N79").Value = -6126

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3025.75
$ws.Range("I80").Value = 3451.3333
$ws.Range("K80").Value = 10353.9999
$ws.Range("M80").Value = -9355.999899999999

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 3025.75
$ws.Range("I83").Value = 3451.3333
$ws.Range("K83").Value = 31061.9997
$ws.Range("M83").Value = -26069.9997

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1334.125
$ws.Range("J88").Value = 1357.6
$ws.Range("L88").Value = 1357.6
$ws.Range("N88").Value = -2169.6

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1334.125
$ws.Range("J91").Value = 1357.6
$ws.Range("L91").Value = 1357.6
$ws.Range("N91").Value = -4165.6

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 426
$ws.Range("I96").Value = 235.14285
$ws.Range("K96").Value = 705.4285500000001
$ws.Range("M96").Value = 667.5714499999999

# ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 2099.2
$ws.Range("J101").Value = 499.66666
$ws.Range("L101").Value = 1498.99998
$ws.Range("N101").Value = -4742.999980000001

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 5136.25
$ws.Range("I107").Value = 3515
$ws.Range("K107").Value = 3515
$ws.Range("M107").Value = -1595

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1201.2
$ws.Range("I132").Value = 1226.8889
$ws.Range("J132").Value = 970
$ws.Range("K132").Value = 3680.6667
$ws.Range("L132").Value = 2910
$ws.Range("M132").Value = -1150.6667
$ws.Range("N132").Value = -7970

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1750
$ws.Range("I137").Value = 1750
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5250
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -2700

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3248.1
$ws.Range("I63").Value = 1386.1111
$ws.Range("J63").Value = 20006
$ws.Range("K63").Value = 1386.1111
$ws.Range("L63").Value = 20006
$ws.Range("M63").Value = -700.1111000000001
$ws.Range("N63").Value = -21378

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3248.1
$ws.Range("I66").Value = 1386.1111
$ws.Range("J66").Value = 20006
$ws.Range("K66").Value = 6930.5555
$ws.Range("L66").Value = 100030
$ws.Range("M66").Value = -3498.5555
$ws.Range("N66").Value = -106894

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1971.0667
$ws.Range("I110").Value = 2124.1667
$ws.Range("J110").Value = 1358.6666
$ws.Range("K110").Value = 2124.1667
$ws.Range("L110").Value = 1358.6666
$ws.Range("M110").Value = -79.16670000000022
$ws.Range("N110").Value = -5448.6666

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3324.875
$ws.Range("I122").Value = 3324.875
$ws.Range("K122").Value = 9974.625
$ws.Range("M122").Value = -7524.625

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3130.818
$ws.Range("J86").Value = 2234.6
$ws.Range("L86").Value = 2234.6
$ws.Range("N86").Value = -4480.6

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3130.818
$ws.Range("J89").Value = 2234.6
$ws.Range("L89").Value = 11173
$ws.Range("N89").Value = -22405

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6601.6
$ws.Range("I105").Value = 6601.6
$ws.Range("K105").Value = 6601.6
$ws.Range("M105").Value = -4854.6

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6360.8125
$ws.Range("I86").Value = 5860.75
$ws.Range("J86").Value = 6860.875
$ws.Range("K86").Value = 5860.75
$ws.Range("L86").Value = 6860.875
$ws.Range("M86").Value = -4737.75
$ws.Range("N86").Value = -9106.875

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 6360.8125
$ws.Range("I89").Value = 5860.75
$ws.Range("J89").Value = 6860.875
$ws.Range("K89").Value = 29303.75
$ws.Range("L89").Value = 34304.375
$ws.Range("M89").Value = -23687.75
$ws.Range("N89").Value = -45536.375

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5985.3335
$ws.Range("J99").Value = 5000
$ws.Range("L99").Value = 5000
$ws.Range("N99").Value = -7996

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 5470.2856
$ws.Range("I105").Value = 5659.4
$ws.Range("K105").Value = 5659.4
$ws.Range("M105").Value = -3912.4

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1721
$ws.Range("I107").Value = 1954.091
$ws.Range("K107").Value = 1954.091
$ws.Range("M107").Value = -34.09099999999989

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5985.3335
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 846
$ws.Range("I68").Value = 846
$ws.Range("K68").Value = 2538
$ws.Range("M68").Value = -1727

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 846
$ws.Range("I71").Value = 846
$ws.Range("K71").Value = 7614
$ws.Range("M71").Value = -3558

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 3111.8333
$ws.Range("I98").Value = 4260.6665
$ws.Range("J98").Value = 1963
$ws.Range("K98").Value = 12781.9995
$ws.Range("L98").Value = 5889
$ws.Range("M98").Value = -11283.9995
$ws.Range("N98").Value = -8885

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 297.66666
$ws.Range("I107").Value = 299.5
$ws.Range("J107").Value = 294
$ws.Range("K107").Value = 898.5
$ws.Range("L107").Value = 882
$ws.Range("M107").Value = 1021.5
$ws.Range("N107").Value = -4722

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3377.75
$ws.Range("I122").Value = 1426
$ws.Range("K122").Value = 12834
$ws.Range("M122").Value = -10384

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1336.125
$ws.Range("J132").Value = 1197.5
$ws.Range("L132").Value = 10777.5
$ws.Range("N132").Value = -15837.5

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5761.4
$ws.Range("I122").Value = 6499.75
$ws.Range("J122").Value = 2808
$ws.Range("K122").Value = 19499.25
$ws.Range("L122").Value = 8424
$ws.Range("M122").Value = -17049.25
$ws.Range("N122").Value = -13324

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3576.2
$ws.Range("I126").Value = 3835.4285
$ws.Range("K126").Value = 11506.2855
$ws.Range("M126").Value = -9036.2855

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2826.7
$ws.Range("I132").Value = 2381.2856
$ws.Range("J132").Value = 3866
$ws.Range("K132").Value = 7143.8568
$ws.Range("L132").Value = 11598
$ws.Range("M132").Value = -4613.8568
$ws.Range("N132").Value = -16658

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3709
$ws.Range("I82").Value = 4380.25
$ws.Range("K82").Value = 4380.25
$ws.Range("M82").Value = -4019.25

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3709
$ws.Range("I85").Value = 4380.25
$ws.Range("K85").Value = 4380.25
$ws.Range("M85").Value = -3132.25

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3757.6
$ws.Range("I93").Value = 3757.6
$ws.Range("K93").Value = 3757.6
$ws.Range("M93").Value = -2509.6

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4370.8887
$ws.Range("I132").Value = 3557
$ws.Range("J132").Value = 5998.6665
$ws.Range("K132").Value = 10671
$ws.Range("L132").Value = 17995.9995
$ws.Range("M132").Value = -8141
$ws.Range("N132").Value = -23055.9995

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3665.5
$ws.Range("I136").Value = 3997
$ws.Range("J136").Value = 3499.75
$ws.Range("K136").Value = 11991
$ws.Range("L136").Value = 10499.25
$ws.Range("M136").Value = -9441
$ws.Range("N136").Value = -15599.25
